# Update teams-position.xlsx:
#  - refresh the standings table (rows 5-18) with the latest games/wins/
#    losses/score/points figures after the new round of matches
#  - append the two new game days (16 & 17 Nov) with their match results,
#    including the newly played "Eagles - SIRIUS" game

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Standings table (B4:H18). Columns: B=Место, C=Команда, D=Игры,
# E=Побед, F=Поражений, G=Мячи, H=Очки. Team order in column C is
# unchanged; only the per-team numbers move.
# ---------------------------------------------------------------------

$ws.Range("C5").Value = "GOLDEN HILL"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "285 - 240"
$ws.Range("H5").Value = 8

$ws.Range("C6").Value = "ISsoft"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "235 - 214"
$ws.Range("H6").Value = 7

$ws.Range("C7").Value = "Грушвиль"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "317 - 262"
$ws.Range("H7").Value = 7

$ws.Range("C8").Value = "Эра-Недвижимости плюс"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "293 - 246"
$ws.Range("H8").Value = 7

$ws.Range("C9").Value = "БГУФК"
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "239 - 199"
$ws.Range("H9").Value = 6

$ws.Range("C10").Value = "VSS"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = "246 - 257"
$ws.Range("H10").Value = 6

$ws.Range("C11").Value = "ОПЛАТИ"
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = "258 - 239"
$ws.Range("H11").Value = 6

$ws.Range("C12").Value = "SIRIUS"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = "257 - 234"
$ws.Range("H12").Value = 6

$ws.Range("C13").Value = "Стрела"
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = "254 - 267"
$ws.Range("H13").Value = 5

$ws.Range("C14").Value = "Mapogo males"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "227 - 231"
$ws.Range("H14").Value = 5

$ws.Range("C15").Value = "NORD"
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = "201 - 354"
$ws.Range("H15").Value = 5

$ws.Range("C16").Value = "Eagles"
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = "214 - 222"
$ws.Range("H16").Value = 5

$ws.Range("C17").Value = "ЛФК"
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = "240 - 268"
$ws.Range("H17").Value = 5

$ws.Range("C18").Value = "Минск 7х"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = "149 - 182"
$ws.Range("H18").Value = 3

# ---------------------------------------------------------------------
# Append the two new game days below the existing schedule (row 45 was
# the previous last row). Each row is merged B:H first (matching the
# existing rows above), then the format is re-cloned from the last
# existing day block (rows 41-45: date row style + match-line row
# style) so the merge border doesn't leave stray styling, and finally
# the value is written.
# ---------------------------------------------------------------------

# --- Saturday 16 Nov 2024 (serial 45612) ---
$ws.Range("B46:H46").Merge()
$ws.Range("B41:H41").Copy()
$ws.Range("B46:H46").PasteSpecial(-4122)
$ws.Range("B46").Value = 45612

$ws.Range("B47:H47").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B47:H47").PasteSpecial(-4122)
$ws.Range("B47").Value = "ЛФК - Грушвиль 70:76 (16:30, БНТУ)"

$ws.Range("B48:H48").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B48:H48").PasteSpecial(-4122)
$ws.Range("B48").Value = "ОПЛАТИ - Минск 7х 63:48 (18:00, БНТУ)"

$ws.Range("B49:H49").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B49:H49").PasteSpecial(-4122)
$ws.Range("B49").Value = "VSS - БГУФК 59:53 (19:30, БНТУ)"

# --- Sunday 17 Nov 2024 (serial 45613) ---
$ws.Range("B50:H50").Merge()
$ws.Range("B41:H41").Copy()
$ws.Range("B50:H50").PasteSpecial(-4122)
$ws.Range("B50").Value = 45613

$ws.Range("B51:H51").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B51:H51").PasteSpecial(-4122)
$ws.Range("B51").Value = "NORD - Стрела 66:65 (11:00, БНТУ)"

$ws.Range("B52:H52").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B52:H52").PasteSpecial(-4122)
$ws.Range("B52").Value = "GOLDEN HILL - Mapogo males 86:69 (12:30, БНТУ)"

$ws.Range("B53:H53").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B53:H53").PasteSpecial(-4122)
$ws.Range("B53").Value = "ISsoft - Эра-Недвижимости плюс 73:60 (14:00, БНТУ)"

$ws.Range("B54:H54").Merge()
$ws.Range("B42:H42").Copy()
$ws.Range("B54:H54").PasteSpecial(-4122)
$ws.Range("B54").Value = "Eagles - SIRIUS 42:48 (15:30, БНТУ)"
